# Update weekly price data for "Hortaliza, Feria Lagunitas de Puerto Montt - Coliflor"
# Existing rows 342-353 are shifted down by two rows (to 344-355), and two new rows of
# fresher data are inserted at 342-343 (a new weekly report).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift existing rows 342..353 down to 344..355 (process bottom-up so we
#     never overwrite a source row before it has been read). Only columns D, I, J, K,
#     L, M, P vary row to row for this block; all the other columns hold identical
#     values (mercado/categoria/etc.) for every row, so they do not need copying.
for ($src = 353; $src -ge 342; $src--) {
    $dst = $src + 2

    $d = $ws.Cells.Item($src, 4).Value2
    $i = $ws.Cells.Item($src, 9).Value2
    $j = $ws.Cells.Item($src, 10).Value2
    $k = $ws.Cells.Item($src, 11).Value2
    $l = $ws.Cells.Item($src, 12).Value2
    $m = $ws.Cells.Item($src, 13).Value2
    $p = $ws.Cells.Item($src, 16).Value2

    $ws.Cells.Item($dst, 1).Value2 = 4
    $ws.Cells.Item($dst, 2).Value2 = "Feria Lagunitas de Puerto Montt"
    $ws.Cells.Item($dst, 3).Value2 = "Los Lagos"
    $ws.Cells.Item($dst, 4).Value2 = $d
    $ws.Cells.Item($dst, 5).Value2 = 10
    $ws.Cells.Item($dst, 6).Value2 = 100112008
    $ws.Cells.Item($dst, 7).Value2 = "Coliflor"
    $ws.Cells.Item($dst, 8).Value2 = "Sin especificar"
    $ws.Cells.Item($dst, 9).Value2 = $i
    $ws.Cells.Item($dst, 10).Value2 = $j
    $ws.Cells.Item($dst, 11).Value2 = $k
    $ws.Cells.Item($dst, 12).Value2 = $l
    $ws.Cells.Item($dst, 13).Value2 = $m
    $ws.Cells.Item($dst, 14).Value2 = "`$/unidad"
    $ws.Cells.Item($dst, 15).Value2 = "Región Metropolitana"
    $ws.Cells.Item($dst, 16).Value2 = $p
    $ws.Cells.Item($dst, 17).Value2 = 1
    $ws.Cells.Item($dst, 18).Value2 = "Hortaliza"

    # Column D keeps the date number format used throughout the column.
    $ws.Cells.Item($dst, 4).NumberFormat = $ws.Cells.Item($src, 4).NumberFormat
}

# --- Step 2: overwrite rows 342 and 343 with the new weekly report entries.
$ws.Cells.Item(342, 1).Value2 = 4
$ws.Cells.Item(342, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(342, 3).Value2 = "Los Lagos"
$ws.Cells.Item(342, 4).Value2 = 44747
$ws.Cells.Item(342, 5).Value2 = 10
$ws.Cells.Item(342, 6).Value2 = 100112008
$ws.Cells.Item(342, 7).Value2 = "Coliflor"
$ws.Cells.Item(342, 8).Value2 = "Sin especificar"
$ws.Cells.Item(342, 9).Value2 = "Primera"
$ws.Cells.Item(342, 10).Value2 = 600
$ws.Cells.Item(342, 11).Value2 = 1800
$ws.Cells.Item(342, 12).Value2 = 1800
$ws.Cells.Item(342, 13).Value2 = 1800
$ws.Cells.Item(342, 14).Value2 = "`$/unidad"
$ws.Cells.Item(342, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(342, 16).Value2 = 1800
$ws.Cells.Item(342, 17).Value2 = 1
$ws.Cells.Item(342, 18).Value2 = "Hortaliza"

$ws.Cells.Item(343, 1).Value2 = 4
$ws.Cells.Item(343, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(343, 3).Value2 = "Los Lagos"
$ws.Cells.Item(343, 4).Value2 = 44747
$ws.Cells.Item(343, 5).Value2 = 10
$ws.Cells.Item(343, 6).Value2 = 100112008
$ws.Cells.Item(343, 7).Value2 = "Coliflor"
$ws.Cells.Item(343, 8).Value2 = "Sin especificar"
$ws.Cells.Item(343, 9).Value2 = "Segunda"
$ws.Cells.Item(343, 10).Value2 = 600
$ws.Cells.Item(343, 11).Value2 = 1500
$ws.Cells.Item(343, 12).Value2 = 1500
$ws.Cells.Item(343, 13).Value2 = 1500
$ws.Cells.Item(343, 14).Value2 = "`$/unidad"
$ws.Cells.Item(343, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(343, 16).Value2 = 1500
$ws.Cells.Item(343, 17).Value2 = 1
$ws.Cells.Item(343, 18).Value2 = "Hortaliza"

# Apply the same date number format to the two new date cells as the rest of column D.
$ws.Cells.Item(342, 4).NumberFormat = $ws.Cells.Item(344, 4).NumberFormat
$ws.Cells.Item(343, 4).NumberFormat = $ws.Cells.Item(344, 4).NumberFormat
